$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new record as row 50, pushing the existing row 50 (and
# everything below it) down by one row.
$ws.Rows.Item(50).Insert()

# Populate the new row 50 with the new weekly record.
$ws.Range("A50").Value = 1
$ws.Range("B50").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C50").Value = "Arica y Parinacota"
$ws.Range("D50").Value = 44895
$ws.Range("E50").Value = 15
$ws.Range("F50").Value = 100112008
$ws.Range("G50").Value = "Coliflor"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Tercera"
$ws.Range("J50").Value = 1500
$ws.Range("K50").Value = 350
$ws.Range("L50").Value = 450
$ws.Range("M50").Value = 400
$ws.Range("N50").Value = "$/unidad"
$ws.Range("O50").Value = "Región de Arica y Parinacota"
$ws.Range("P50").Value = 400
$ws.Range("Q50").Value = 1
$ws.Range("R50").Value = "Hortaliza"
